# Update odds values for the week-5 and week-6 match rows in the FlashScore sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 updates
$ws.Range("G5").Value = 1.5
$ws.Range("H5").Value = 3.75
$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 2.1
$ws.Range("L5").Value = 7
$ws.Range("U5").Value = 2.25
$ws.Range("V5").Value = 1.57
$ws.Range("Z5").Value = 10
$ws.Range("AD5").Value = 7.5
$ws.Range("AE5").Value = 21
$ws.Range("AH5").Value = 34
$ws.Range("AJ5").Value = 81
$ws.Range("AK5").Value = 51
$ws.Range("AN5").Value = 3.25
$ws.Range("AO5").Value = 8
$ws.Range("AQ5").Value = 26
$ws.Range("AW5").Value = 8
$ws.Range("AZ5").Value = 151
$ws.Range("BA5").Value = 201

# Row 6 updates
$ws.Range("AC6").Value = 8
$ws.Range("AI6").Value = 17
$ws.Range("AL6").Value = 51
$ws.Range("AY6").Value = 41
